$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.474.56"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "'3.023.83"
$ws.Range("E3").Value = "  +4.22%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'200.18"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'627.60"
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.553"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'0.208"
$ws.Range("E9").Value = "  +4.96%  "
$ws.Range("D10").Value = "'3.024.52"
$ws.Range("E10").Value = "  +4.15%  "
$ws.Range("D11").Value = "'0.435"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'5.15"
$ws.Range("E13").Value = "  +6.10%  "
$ws.Range("D14").Value = "'3.578.38"
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").Value = "'29.31"
$ws.Range("E15").Value = "  +6.89%  "
$ws.Range("D16").Value = "'76.338.06"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'3.024.73"
$ws.Range("E18").Value = "  +4.21%  "
$ws.Range("D19").Value = "'13.47"
$ws.Range("E19").Value = "  +5.79%  "
$ws.Range("D20").Value = "'8.98"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "'373.69"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'4.35"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "'73.15"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").Value = "'3.181.09"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +4.97%  "
$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("D29").Value = "'0.0000110"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "'506.71"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  +7.46%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'20.71"
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").Value = "'163.98"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "'0.384"
$ws.Range("E38").Value = "  +11.39%  "
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").Value = "'189.67"
$ws.Range("E40").Value = "  +5.36%  "
$ws.Range("D41").Value = "'0.105"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'5.12"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").Value = "'42.40"
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("E47").Value = "  +6.16%  "
$ws.Range("D48").Value = "'0.715"
$ws.Range("E48").Value = "  +9.22%  "
$ws.Range("D49").Value = "'0.602"
$ws.Range("E49").Value = "  +5.76%  "
$ws.Range("D50").Value = "'2.38"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("E51").Value = "  +4.64%  "
